$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 58, shifting the existing
# rows 58:75 down to 59:76.
$ws.Rows.Item(58).Insert()

# Populate the new row 58 with the new weekly price record.
$ws.Range("A58").Value = 5
$ws.Range("B58").Value = "Macroferia Regional de Talca"
$ws.Range("C58").Value = "Maule"
$ws.Range("D58").Value = 44588
$ws.Range("E58").Value = 7
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = "Frutos de hueso (carozo)"
$ws.Range("I58").Value = 100103002
$ws.Range("J58").Value = "Ciruela"
$ws.Range("K58").Value = "Black Amber"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 200
$ws.Range("N58").Value = 9000
$ws.Range("O58").Value = 9000
$ws.Range("P58").Value = 9000
$ws.Range("Q58").Value = "$/bandeja 18 kilos granel"
$ws.Range("R58").Value = "Provincia de Curicó"
$ws.Range("S58").Value = 500
$ws.Range("T58").Value = 18
